# Applies the cached-value refresh captured in the commit diff.
# These H:N columns hold plain numeric <v> caches (no formulas) that were
# regenerated by the scheduled market-data runner, so we just overwrite the
# literal cell values (and clear the handful of cells that became blank, or
# add the few that became populated).
$wb = $excel.ActiveWorkbook

$sheetChanges = @{
    "ALC" = @{
        "H19" = 2733.6667
        "J19" = 2300
        "L19" = 2300
        "N19" = -2650
        "H32" = 8801.375
        "I32" = 6500
        "K32" = 6500
        "M32" = -6174
        "H40" = 999
        "I40" = 999
        "K40" = 999
        "M40" = -824
        "H76" = 1209.25
        "I76" = 1282.3334
        "J76" = 990
        "K76" = 1282.3334
        "L76" = 990
        "M76" = -967.3334
        "N76" = -1620
        "H79" = 1209.25
        "I79" = 1282.3334
        "J79" = 990
        "K79" = 1282.3334
        "L79" = 990
        "M79" = -190.3334
        "N79" = -3174
        "H112" = 1990.75
        "J112" = 1990.75
        "L112" = 5972.25
        "N112" = -8188.25
        "H119" = 0
        "I119" = 0
        "J119" = 0
        "K119" = 0
        "L119" = 0
        "M119" = $null
        "N119" = $null
        "H137" = 4049.1667
        "I137" = 2824.5
        "J137" = 6498.5
        "K137" = 8473.5
        "L137" = 19495.5
        "M137" = -5923.5
        "N137" = -24595.5
    }
    "ARM" = @{
        "H61" = 2002.75
        "I61" = 1670.3334
        "K61" = 1670.3334
        "M61" = -1458.3334
        "H74" = 2577.0698
        "I74" = 2275.9443
        "K74" = 2275.9443
        "M74" = -1401.9443
        "H77" = 2577.0698
        "I77" = 2275.9443
        "K77" = 11379.7215
        "M77" = -7011.7215
        "H132" = 4666.0835
        "I132" = 0
        "J132" = 4666.0835
        "K132" = 0
        "L132" = 13998.2505
        "M132" = $null
        "N132" = -19058.2505
        "H134" = 0
        "J134" = 0
        "L134" = 0
        "N134" = $null
        "H136" = 2002.75
        "I136" = 1670.3334
        "K136" = 5011.0002
        "M136" = -2461.0002
        "H137" = 100000
        "J137" = 100000
        "L137" = 100000
        "N137" = -110200
    }
    "BSM" = @{
        "H20" = 6932.5557
        "I20" = 6484.857
        "J20" = 8499.5
        "K20" = 6484.857
        "L20" = 8499.5
        "M20" = -6237.857
        "N20" = -8993.5
        "H59" = 99950
        "J59" = 99950
        "L59" = 99950
        "N59" = -101644
        "H105" = 4655.3335
        "I105" = 4655.3335
        "K105" = 4655.3335
        "M105" = -2908.3335
        "H134" = 2678.1936
        "I134" = 2715.6333
        "K134" = 8146.8999
        "M134" = -5611.8999
    }
    "CRP" = @{
        "H31" = 2255.186
        "I31" = 2209.158
        "J31" = 2605
        "K31" = 2209.158
        "L31" = 2605
        "M31" = -1914.158
        "N31" = -3195
        "H34" = 2255.186
        "I34" = 2209.158
        "J34" = 2605
        "K34" = 2209.158
        "L34" = 2605
        "M34" = -2007.158
        "N34" = -3009
        "H122" = 3028.6667
        "I122" = 2709.6667
        "K122" = 8129.000100000001
        "M122" = -5679.000100000001
        "H132" = 3393.6667
        "I132" = 2991.111
        "K132" = 8973.332999999999
        "M132" = -6443.332999999999
    }
    "CUL" = @{
        "H118" = 4564.6665
        "J118" = 5000
        "L118" = 15000
        "N118" = -17486
        "H140" = 2381.5
        "I140" = 2381.5
        "K140" = 7144.5
        "M140" = -1964.5
    }
    "GSM" = @{
        "H70" = 14979
        "I70" = 14998.75
        "J70" = 14900
        "K70" = 14998.75
        "L70" = 14900
        "M70" = -14728.75
        "N70" = -15440
        "H73" = 14979
        "I73" = 14998.75
        "J73" = 14900
        "K73" = 14998.75
        "L73" = 14900
        "M73" = -14062.75
        "N73" = -16772
        "H107" = 567
        "I107" = 567.3333
        "K107" = 567.3333
        "M107" = 1352.6667
    }
    "LTW" = @{
        "H93" = 4054.5
        "I93" = 4054.5
        "J93" = 0
        "K93" = 4054.5
        "L93" = 0
        "M93" = -2806.5
        "N93" = $null
        "H137" = 83000
        "J137" = 83000
        "L137" = 83000
        "N137" = -93200
    }
    "WVR" = @{
        "H126" = 2105
        "I126" = 1909.5
        "K126" = 5728.5
        "M126" = -3258.5
        "H136" = 1534.9375
        "I136" = 1201.3
        "K136" = 3603.9
        "M136" = -1053.9
    }
}

foreach ($sheetName in $sheetChanges.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $sheetChanges[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $newValue = $cellMap[$cellRef]
        if ($null -eq $newValue) {
            $ws.Range($cellRef).ClearContents()
        } else {
            $ws.Range($cellRef).Value = $newValue
        }
    }
}

Write-Host "Applied cell updates across $($sheetChanges.Keys.Count) sheets."
